$d = $word.ActiveDocument
$d.Content.Find.Execute("Criar um sistema destinado a salões de cabeleireiro para gerenciar informações de profissionais, procedimentos, agendamentos a serem realizados e de clientes.", $false, $false, $false, $false, $false, $true, 1, $false, "Criar um sistema ERP (Enterprise Resource Planning) destinado a salões de cabeleireiro para gerenciar informações de profissionais, procedimentos, clientes e de agendamentos a serem realizados.", 2)
Write-Output "done"
